$d = $word.ActiveDocument

$map = [ordered]@{
    "73×42=3066" = "33×67=2211"
    "62×22=1364" = "36×93=3348"
    "25×50=1250" = "97×11=1067"
    "23×30=690"  = "46×24=1104"
    "43×78=3354" = "35×42=1470"
    "91×82=7462" = "15×15=225"
    "95×72=6840" = "97×26=2522"
    "68×38=2584" = "71×37=2627"
    "89×46=4094" = "63×19=1197"
    "47×50=2350" = "32×26=832"
    "89×41=3649" = "83×47=3901"
    "56×81=4536" = "42×29=1218"
    "43×57=2451" = "87×85=7395"
    "44×25=1100" = "82×50=4100"
    "39×67=2613" = "45×81=3645"
    "57×22=1254" = "31×24=744"
    "73×64=4672" = "60×11=660"
    "24×47=1128" = "56×72=4032"
    "98×65=6370" = "35×38=1330"
    "22×23=506"  = "12×38=456"
    "21×62=1302" = "48×50=2400"
    "80×56=4480" = "40×39=1560"
    "39×53=2067" = "86×15=1290"
    "44×51=2244" = "13×20=260"
    "51×38=1938" = "63×37=2331"
}

foreach ($key in $map.Keys) {
    $old = $key
    $new = $map[$key]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
